$d = $word.ActiveDocument

# Locate the paragraph that ends the "Get the starter version of the class."
# bullet item - the new bullet about <summary> is inserted right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Get the starter version of the class.*") {
        $target = $p
    }
}

# Insert a new paragraph right after it. Word automatically carries over the
# paragraph formatting (pStyle/numPr) of the paragraph it was split from.
$target.Range.InsertParagraphAfter() | Out-Null

# Re-fetch the (now empty) paragraph that was just created, it directly
# follows the "Get the starter..." paragraph.
$newPara = $target.Next()

# Build the run content for the new bullet via a WordprocessingML fragment so
# we get full control of the run properties (font, size, color, highlight)
# that are applied to the "<summary>" text, while keeping the paragraph
# properties (pStyle/numPr) that are already on the empty paragraph.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:pPr>' +
'<w:pStyle w:val="ListParagraph"/>' +
'<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
'</w:pPr>' +
'<w:r><w:t xml:space="preserve">Add </w:t></w:r>' +
'<w:r>' +
'<w:rPr>' +
'<w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/>' +
'<w:color w:val="808080"/>' +
'<w:sz w:val="19"/>' +
'<w:szCs w:val="19"/>' +
'<w:highlight w:val="white"/>' +
'</w:rPr>' +
'<w:t>&lt;summary&gt;</w:t>' +
'</w:r>' +
'<w:r><w:t xml:space="preserve"> as good practice.</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xml) | Out-Null
